$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.369.67'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '2.233.85'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.629'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.45'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.61%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.621'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '43.12'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0963'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.64%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.12'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.49%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.103'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.45'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.60%  '
$ws.Range('E15').Value = '  -0.55%  '
$ws.Range('D16').Value = '2.235.55'
$ws.Range('E16').Value = '  -1.25%  '
$ws.Range('D17').Value = '42.179.34'
$ws.Range('E17').Value = '  +0.36%  '
$ws.Range('E18').Value = '  +14.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.16'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.36%  '
$ws.Range('E20').Value = '  +0.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.24'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +39.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '230.93'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('E23').Value = '  -5.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.75'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.15%  '
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.67'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.30'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.78%  '
$ws.Range('E28').Value = '  +6.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '166.47'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.94'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.82'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +17.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0807'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.118'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.92%  '
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.76'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -8.09%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.125'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.39'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.64%  '
$ws.Range('E37').Value = '  +2.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '13.28'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.72%  '
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.63'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '63.26'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.201'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.71%  '
$ws.Range('E43').Value = '  +2.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '104.89'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.67%  '
$ws.Range('E45').Value = '  +3.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.995'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.17%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.38'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.53%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.13'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.33%  '
$ws.Range('E49').Value = '  +0.72%  '
$ws.Range('E50').Value = '  +0.73%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.06'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.24%  '
